$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0025387765691307
$ws.Range("C2").Value = 0.996960880603376
$ws.Range("D2").Value = 0.00657858162074014
$ws.Range("E2").Value = 0.0000741248633322832
$ws.Range("F2").Value = 0.000611530122491337
$ws.Range("G2").Value = 0.000259437021662991
$ws.Range("H2").Value = 0.000166780942497637
$ws.Range("I2").Value = 0.00159368456164409
$ws.Range("J2").Value = 0.00203843374163779
$ws.Range("K2").Value = 0.00124159146081574
$ws.Range("L2").Value = 0.997479754646702
$ws.Range("M2").Value = 0.94423957155829
$ws.Range("N2").Value = 0.0273891370012787
$ws.Range("O2").Value = 0.0000185312158330708
$ws.Range("P2").Value = 0.0119896966439968
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.997461223430869
$ws.Range("S2").Value = 0.940051516780016
$ws.Range("T2").Value = 0.0000370624316661416
$ws.Range("U2").Value = 0.0000185312158330708
$ws.Range("V2").Value = 0.978133165316976
$ws.Range("W2").Value = 0.00281674480662676
$ws.Range("X2").Value = 0.00333561884995275
$ws.Range("B3").Value = 0.000111187294998425
$ws.Range("C3").Value = 0.000574467690825195
$ws.Range("D3").Value = 0.892926634916517
$ws.Range("E3").Value = 0.00092656079165354
$ws.Range("F3").Value = 0.000407686748327558
$ws.Range("G3").Value = 0.883623964568315
$ws.Range("H3").Value = 0.00187165279914015
$ws.Range("I3").Value = 0.00148249726664566
$ws.Range("J3").Value = 0.000352093100828345
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.00137130997164724
$ws.Range("M3").Value = 0.0000741248633322832
$ws.Range("N3").Value = 0.0000555936474992124
$ws.Range("O3").Value = 0.00046328039582677
$ws.Range("P3").Value = 0.0657116913440691
$ws.Range("Q3").Value = 0.999796156625836
$ws.Range("R3").Value = 0.000592998906658266
$ws.Range("S3").Value = 0.0000555936474992124
$ws.Range("T3").Value = 0.0417322980560755
$ws.Range("U3").Value = 0.000092656079165354
$ws.Range("V3").Value = 0.0000555936474992124
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.0000370624316661416
$ws.Range("B4").Value = 0.996942349387543
$ws.Range("C4").Value = 0.00100068565498582
$ws.Range("D4").Value = 0.00135277875581417
$ws.Range("E4").Value = 0.000166780942497637
$ws.Range("F4").Value = 0.0000741248633322832
$ws.Range("G4").Value = 0.000129718510831496
$ws.Range("H4").Value = 0.000611530122491337
$ws.Range("I4").Value = 0.995737820358394
$ws.Range("J4").Value = 0.99179067138595
$ws.Range("K4").Value = 0.99673850601338
$ws.Range("L4").Value = 0.000630061338324407
$ws.Range("M4").Value = 0.0474955061801605
$ws.Range("N4").Value = 0.966532624205474
$ws.Range("O4").Value = 0.0000741248633322832
$ws.Range("P4").Value = 0.0255360154179716
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.000518874043325983
$ws.Range("S4").Value = 0.0525174656709227
$ws.Range("T4").Value = 0.0000370624316661416
$ws.Range("U4").Value = 0.0000185312158330708
$ws.Range("V4").Value = 0.0186609343439023
$ws.Range("W4").Value = 0.997072067898375
$ws.Range("X4").Value = 0.996479068991717
$ws.Range("B5").Value = 0.0000185312158330708
$ws.Range("C5").Value = 0.00137130997164724
$ws.Range("D5").Value = 0.0845764690621352
$ws.Range("E5").Value = 0.998721346107518
$ws.Range("F5").Value = 0.99879547097085
$ws.Range("G5").Value = 0.100921001426904
$ws.Range("H5").Value = 0.997220317625039
$ws.Range("I5").Value = 0.000833904712488186
$ws.Range("J5").Value = 0.00479958490076534
$ws.Range("K5").Value = 0.00189018401497322
$ws.Range("L5").Value = 0.000407686748327558
$ws.Range("M5").Value = 0.0000555936474992124
$ws.Range("N5").Value = 0.000092656079165354
$ws.Range("O5").Value = 0.999407001093342
$ws.Range("P5").Value = 0.882122936085837
$ws.Range("Q5").Value = 0.000092656079165354
$ws.Range("R5").Value = 0.00135277875581417
$ws.Range("S5").Value = 0.0000370624316661416
$ws.Range("T5").Value = 0.953190148805663
$ws.Range("U5").Value = 0.999796156625836
$ws.Range("V5").Value = 0.000092656079165354
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.0000185312158330708
